$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 403, shifting existing rows 403:484 down to 404:485.
$ws.Rows("403").Insert()

# Populate the newly inserted row 403 with the new data record.
$ws.Range("A403").Value2 = 5
$ws.Range("B403").Value2 = "Macroferia Regional de Talca"
$ws.Range("C403").Value2 = "Maule"
$ws.Range("D403").Value2 = 44694
$ws.Range("E403").Value2 = 7
$ws.Range("F403").Value2 = 100114001
$ws.Range("G403").Value2 = "Papa"
$ws.Range("H403").Value2 = "Patagonia"
$ws.Range("I403").Value2 = "1a (cosecha)"
$ws.Range("J403").Value2 = 1500
$ws.Range("K403").Value2 = 8000
$ws.Range("L403").Value2 = 8000
$ws.Range("M403").Value2 = 8000
$ws.Range("N403").Value2 = "`$/saco 25 kilos"
$ws.Range("O403").Value2 = "Región de Los Lagos"
$ws.Range("P403").Value2 = 320
$ws.Range("Q403").Value2 = 25
$ws.Range("R403").Value2 = "Hortaliza"
